# Update "想去人数" (want-to-go count) figures that changed between scrapes.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 10871
$ws1.Range("F4").Value = 70
$ws1.Range("F5").Value = 740
$ws1.Range("F6").Value = 503

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 10871
$ws4.Range("F4").Value = 70
$ws4.Range("F5").Value = 740
$ws4.Range("F7").Value = 503
